# Add the 2024/11/09 data column (BJ) to the "合成確率" sheet,
# mirroring the existing per-day columns (header + 52 data rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# PasteSpecial / paste-type constants (standard Excel enum values).
$xlPasteFormats = -4122
$xlPasteValues  = -4163

# New column (BJ = 62nd column): same width as the previous date column (BI).
$ws.Columns.Item(62).ColumnWidth = $ws.Columns.Item(61).ColumnWidth

# Header cell BJ1 holds the date as text, styled like the other date headers (BI1).
# Build it as a text formula first (so "2024/11/09" is never auto-parsed into a real
# date), then flatten the formula down to a plain cached value.
$ws.Range("BI1").Copy()
$ws.Range("BJ1").PasteSpecial($xlPasteFormats)
$ws.Range("BJ1").Formula = "=""2024/11/09"""
$ws.Range("BJ1").Copy()
$ws.Range("BJ1").PasteSpecial($xlPasteValues)

# The sheet uses three pre-existing "highlight" styles for the numeric data
# (depending on the value range): plain / yellow / light-blue fill. Reuse
# existing cells that already carry each style as format templates.
$styleTemplates = @{
    1 = $ws.Range("A2")
    2 = $ws.Range("D2")
    3 = $ws.Range("N2")
}

# New data for 2024/11/09: row number, style id (1/2/3), value.
$newData = @(
    @{ Row = 2; Style = 3; Value = 138.2 }
    @{ Row = 3; Style = 1; Value = 165.9 }
    @{ Row = 4; Style = 2; Value = 122.1 }
    @{ Row = 5; Style = 1; Value = 152.4 }
    @{ Row = 6; Style = 1; Value = 155.3 }
    @{ Row = 7; Style = 1; Value = 160.2 }
    @{ Row = 8; Style = 1; Value = 281.5 }
    @{ Row = 9; Style = 1; Value = 172.3 }
    @{ Row = 10; Style = 1; Value = 190.2 }
    @{ Row = 11; Style = 2; Value = 123.4 }
    @{ Row = 12; Style = 3; Value = 137.1 }
    @{ Row = 13; Style = 3; Value = 126 }
    @{ Row = 14; Style = 3; Value = 129.1 }
    @{ Row = 15; Style = 1; Value = 171.4 }
    @{ Row = 16; Style = 3; Value = 135.7 }
    @{ Row = 17; Style = 3; Value = 138.3 }
    @{ Row = 18; Style = 1; Value = 144.5 }
    @{ Row = 19; Style = 1; Value = 284 }
    @{ Row = 20; Style = 1; Value = 217.7 }
    @{ Row = 21; Style = 1; Value = 176.3 }
    @{ Row = 22; Style = 1; Value = 140.8 }
    @{ Row = 23; Style = 1; Value = 156.5 }
    @{ Row = 24; Style = 1; Value = 230.1 }
    @{ Row = 25; Style = 1; Value = 180.6 }
    @{ Row = 26; Style = 1; Value = 143.5 }
    @{ Row = 27; Style = 1; Value = 165 }
    @{ Row = 28; Style = 1; Value = 141.9 }
    @{ Row = 29; Style = 1; Value = 174 }
    @{ Row = 30; Style = 1; Value = 161.4 }
    @{ Row = 31; Style = 1; Value = 162.1 }
    @{ Row = 32; Style = 1; Value = 204 }
    @{ Row = 33; Style = 1; Value = 151 }
    @{ Row = 34; Style = 1; Value = 153.7 }
    @{ Row = 35; Style = 3; Value = 129.1 }
    @{ Row = 36; Style = 1; Value = 202.1 }
    @{ Row = 37; Style = 2; Value = 113.8 }
    @{ Row = 38; Style = 1; Value = 214 }
    @{ Row = 39; Style = 1; Value = 174.6 }
    @{ Row = 40; Style = 1; Value = 159.9 }
    @{ Row = 41; Style = 1; Value = 142 }
    @{ Row = 42; Style = 1; Value = 172.5 }
    @{ Row = 43; Style = 1; Value = 192.3 }
    @{ Row = 44; Style = 1; Value = 195.6 }
    @{ Row = 45; Style = 1; Value = 141.7 }
    @{ Row = 46; Style = 1; Value = 173.3 }
    @{ Row = 47; Style = 3; Value = 135.1 }
    @{ Row = 48; Style = 1; Value = 144.4 }
    @{ Row = 49; Style = 1; Value = 175.4 }
    @{ Row = 50; Style = 2; Value = 113.2 }
    @{ Row = 51; Style = 2; Value = 122.9 }
    @{ Row = 52; Style = 3; Value = 137.1 }
    @{ Row = 53; Style = 1; Value = 144.1 }
)

foreach ($item in $newData) {
    $target = $ws.Range("BJ" + $item.Row)
    $styleTemplates[$item.Style].Copy()
    $target.PasteSpecial($xlPasteFormats)
    $target.Value = $item.Value
}
